$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "257.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.04%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.96"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.38%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.675"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-10.15%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05880"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.61%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.626"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.80%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8588"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.62%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9486"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-6.85%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1409"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.25%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.03983"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "12.45%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07092"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.29%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03180"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.30%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09159"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.82%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001544"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.39%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006218"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "4.58%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.524"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.78%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.209"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.80%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.206"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.65%"
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.01051"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1,634.57%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3052"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.88%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.48%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.822"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.28%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04230"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.41%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001220"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.08%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004297"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.71%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.05%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001937"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "30.53%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.55%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006200"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "9.79%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1101"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.26%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002200"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.81%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.51%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005446"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.34%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.05%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07000"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-35.80%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2309"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10,498.78%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.05%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.05%"
